$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: "Speed V/S CPU" label and self-ratio formula for B13
$ws.Cells.Item(13, 1).Value = "Speed V/S CPU"
$ws.Cells.Item(13, 2).Formula = "=B12/B12"

# Selection moved from G15 to B14
$ws.Range("B14").Select()

Write-Output "done"
